# Apply edits described in the commit diff to the BBNPPTY sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BBNPPTY")

# Columns B, C and D (years 2021-2023) for rows 2 through 25 all become a
# Boolean "1" (ban new power plants flag) and lose their integer ("0")
# number-format style, reverting to the default/normal style.
$rng = $ws.Range("B2:D25")
$rng.Value = 1
$rng.Style = "Normal"

# Update the sheet's active selection to match the saved view state.
$ws.Range("B2:D25").Select()
